$d = $word.ActiveDocument
$d.Content.Find.Execute("Epic Systems", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Hudson River Trading", 2)
